$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.558.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.392.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.67"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.969.51"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.393.97"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.584.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.00"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.79"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.36%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.90"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "169.25"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.425.58"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.78"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.64"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.462.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0262"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.26%  "
